$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $rng = $cell.Range
    # Exclude the trailing cell-mark character so only the visible text is replaced,
    # preserving the run's formatting (rPr).
    $textRng = $d.Range($rng.Start, $rng.End - 1)
    $textRng.Text = $newText
}

Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "305"
Set-CellText $t 5  "0.00001"
Set-CellText $t 6  "0.00063"
Set-CellText $t 7  "0.00018"
Set-CellText $t 8  "0.00005"
Set-CellText $t 9  "0.00031"
Set-CellText $t 10 "0.00038"
Set-CellText $t 11 "0.00042"
Set-CellText $t 12 "0.06137"

Set-CellText $t 44 "99.7"
Set-CellText $t 45 "0.06"
Set-CellText $t 46 "20"
